{"js": "// Locate the \"page views\" levels table (3 columns: Level / Description / Number of time series)\n// that is immediately followed by the \"Total number of levels = 5\" / \"Total number of time\n// series = 950\" paragraphs, and rewrite it to reflect the new base HF model results:\n//   - Level 1 becomes the old \"Access\" row (count 3 instead of 4)\n//   - Level 2 becomes the old \"Agent\" row (count simplified to \"5 (mobile app doesn't have spider)\")\n//   - Level 3 becomes the old \"Language\" row (count becomes the formula \"4*5 = 20\")\n//   - A brand-new Level 4 \"Purpose\" row (count 153) is inserted\n//   - The old \"Article\" row becomes Level 5 and its count drops the per-item breakdown\n//   - The summary paragraphs are updated: levels 5 -> 6, time series 950 -> 1095\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const t of tables.items) {\n  t.load(\"values\");\n}\nawait context.sync();\n\nlet targetTable = null;\nfor (const t of tables.items) {\n  const vals = t.values;\n  if (\n    vals.length > 0 &&\n    vals[0][0] === \"Level\" &&\n    vals[0][1] === \"Description\" &&\n    vals.some((row) => row[1] && row[1].indexOf(\"Total Page Views\") !== -1)\n  ) {\n    targetTable = t;\n    break;\n  }\n}\n\nif (!targetTable) {\n  throw new Error(\"Could not find the page-views levels table\");\n}\n\nconst rows = targetTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row indices in the existing table (0 = header):\n// 0 header, 1 = Level 0 (Total Page Views), 2 = Level 1 (Language),\n// 3 = Level 2 (Access), 4 = Level 3 (Agent), 5 = Level 4 (Article)\nconst rowLanguage = rows.items[2];\nconst rowAccess = rows.items[3];\nconst rowAgent = rows.items[4];\nconst rowArticle = rows.items[5];\n\nfor (const r of [rowLanguage, rowAccess, rowAgent, rowArticle]) {\n  r.cells.load(\"items\");\n}\nawait context.sync();\n\nfunction setRow(row, level, description, count) {\n  const cells = row.cells.items;\n  cells[0].value = level;\n  cells[1].value = description;\n  cells[2].value = count;\n}\n\n// Rewrite the existing rows in place with the new (reordered) content.\nsetRow(rowLanguage, \"1\", \"Access - desktop, mobileapp, mobile web\", \"3\");\nsetRow(rowAccess, \"2\", \"Agent - spider, user\", \"5 (mobile app doesn\\u2019t have spider)\");\nsetRow(rowAgent, \"3\", \"Language - de, en, es, zh\", \"4*5 = 20\");\nsetRow(rowArticle, \"5\", \"Article\", \"913\");\nawait context.sync();\n\n// Insert the new \"Purpose\" row right before the (now Level 5) Article row.\nrowArticle.insertRows(\"Before\", 1, [[\"4\", \"Purpose\", \"153\"]]);\nawait context.sync();\n\n// Update the summary paragraphs that immediately follow the table.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet levelsPara = null;\nlet seriesPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"Total number of levels = 5\") {\n    levelsPara = p;\n  } else if (p.text === \"Total number of time series = 950\") {\n    seriesPara = p;\n  }\n}\n\nif (levelsPara) {\n  levelsPara.insertText(\"Total number of levels = 6\", \"Replace\");\n}\nif (seriesPara) {\n  seriesPara.insertText(\"Total number of time series = 1095\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Rewrite the \"page views\" levels table (Level / Description / Number of time series)\n# to reflect the new base HF model results:\n#   - Level 1 becomes the old \"Access\" row (count 3 instead of 4)\n#   - Level 2 becomes the old \"Agent\" row (count simplified to \"5 (mobile app doesn't have spider)\")\n#   - Level 3 becomes the old \"Language\" row (count becomes the formula \"4*5 = 20\")\n#   - A brand-new Level 4 \"Purpose\" row (count 153) is inserted\n#   - The old \"Article\" row becomes Level 5 and its count drops the per-item breakdown\n#   - The summary paragraphs are updated: levels 5 -> 6, time series 950 -> 1095\n\n$doc = $word.ActiveDocument\n\n# Find the table whose row 2 reads \"Total Page Views\" - that's the one we need.\n$targetTable = $null\nfor ($i = 1; $i -le $doc.Tables.Count; $i++) {\n    $candidate = $doc.Tables.Item($i)\n    if ($candidate.Cell(2, 2).Range.Text.Contains(\"Total Page Views\")) {\n        $targetTable = $candidate\n        break\n    }\n}\n\n# Rows (1-based) before the edit:\n#   1 header, 2 Level 0, 3 Level 1 (Language), 4 Level 2 (Access), 5 Level 3 (Agent), 6 Level 4 (Article)\n# Insert the four reordered/new rows right before the old \"Language\" row (row 3), then delete\n# the three old rows (Language, Access, Agent) that got pushed further down.\n$anchor = $targetTable.Rows.Item(3)\n\n$newData = @(\n    @(\"1\", \"Access - desktop, mobileapp, mobile web\", \"3\"),\n    @(\"2\", \"Agent - spider, user\", \"5 (mobile app doesn\" + [char]8217 + \"t have spider)\"),\n    @(\"3\", \"Language - de, en, es, zh\", \"4*5 = 20\"),\n    @(\"4\", \"Purpose\", \"153\")\n)\n\nfor ($i = $newData.Length - 1; $i -ge 0; $i--) {\n    $rowVals = $newData[$i]\n    $newRow = $targetTable.Rows.Add($anchor)\n    $newRow.Cells.Item(1).Range.Text = $rowVals[0]\n    $newRow.Cells.Item(2).Range.Text = $rowVals[1]\n    $newRow.Cells.Item(3).Range.Text = $rowVals[2]\n}\n\n# The old Language/Access/Agent rows are now rows 7, 8, 9 - remove them.\nfor ($i = 1; $i -le 3; $i++) {\n    $targetTable.Rows.Item(7).Delete()\n}\n\n# The old \"Article\" row is now row 7. Its description stays \"Article\"; update Level -> 5 and\n# simplify the time-series count to drop the per-item breakdown.\n$articleRow = $targetTable.Rows.Item(7)\n$articleRow.Cells.Item(1).Range.Text = \"5\"\n\n$countCell = $articleRow.Cells.Item(3)\n$countFind = $countCell.Range.Find\n$countFind.Text = \"913 (38 38 30 24 38 93 93 74 92 93 38 38 25 31 38 29 29 14 29 29)\"\n$countFind.Replacement.Text = \"913\"\n$countFind.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1)\n\n# Update the summary paragraphs that immediately follow the table.\n$tblEnd = $targetTable.Range.End\n$afterTable = $doc.Range($tblEnd, $doc.Content.End)\n$find = $afterTable.Find\n$find.Text = \"Total number of levels = 5\"\n$find.Replacement.Text = \"Total number of levels = 6\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1)\n\n$afterTable2 = $doc.Range($tblEnd, $doc.Content.End)\n$find2 = $afterTable2.Find\n$find2.Text = \"Total number of time series = 950\"\n$find2.Replacement.Text = \"Total number of time series = 1095\"\n$find2.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1)\n"}
